# STS IR Bot Dispatcher / Data / Config.xlsx
# - Add a new "PathCustomerNameList" setting (Customer Account List path) as a
#   new row on the Settings sheet.
# - Fix the hard-coded Windows username baked into the "PathTempDirectory"
#   value by replacing it with a generic "<USERNAME>" placeholder.
# - Restore the view/selection state (active sheet/tab + selected cells) that
#   Excel recorded when the file was last saved.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("Constants")

# Insert a brand-new row above the old "PathTempDirectory" row (row 12) so the
# new setting becomes the first one listed in that block; everything below it
# (PathTempDirectory, the blank separator, ReturnsFile_WorksheetName,
# TempFile_FileName, ...) shifts down by one row.
$ws1.Rows.Item(12).Insert() | Out-Null
# Match the 14.25pt row height used by the surrounding rows in this block.
$ws1.Rows.Item(12).RowHeight = 14.25

# Fix the hardcoded local username baked into the temp-directory path (now
# living in row 13 after the insert above).
$ws1.Range("B13").Value = "C:\Users\<USERNAME>\Documents\UiPath\temp"

# Populate the new row with the Customer Account List setting.
$ws1.Range("B12").Value = "Data\Customer Account List.xlsx"
$ws1.Range("A12").Value = "PathCustomerNameList"

# Restore sheet selections: the Constants sheet remembers A15 was last
# selected (no longer the active tab, no longer scrolled away from A1),
# while Settings becomes the active tab with A12 selected.
$ws2.Select() | Out-Null
$ws2.Range("A15").Select() | Out-Null
$ws1.Select() | Out-Null
$ws1.Range("A12").Select() | Out-Null
